# ---------------------------------------------------------------------------
# Edit summary (per the commit's diff):
#   1. In "...roughly 2000-6000. However, upon sampling additional data,
#      Figure 2 illustrates..." change "However" -> "Moreover", with the
#      surrounding run split into three runs: ". " / "Moreover" / ", ".
#   2. Merge the "(Figure 2)" run and the following " " run (which sit next
#      to each other in "...graph above (Figure 2) supports...") into a
#      single run "(Figure 2) ".
# Both edits live in the same (very long) paragraph, so we apply them from
# the rightmost position to the leftmost to avoid later edits clobbering
# earlier ones (this engine re-flows every run from the edit point to the
# end of the paragraph whenever the paragraph's text is touched).
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- locate the unique anchor text for edit #1 (near "...2000-6000.") -----
$editRange = $d.Content
$editRange.Find.Execute("2000-6000. However, upon sampling", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if (-not $editRange.Find.Found) {
    throw "Could not locate the '2000-6000. However, upon sampling' anchor text"
}

# Replace just the sentence punctuation + 'However, ' text with '. Moreover, '
# (scoped to this exact occurrence so the other 'However, ' runs earlier in
# the paragraph are left untouched).
$howeverRange = $d.Range($editRange.Start, $editRange.Start + [int]"2000-6000. However, ".Length)
$howeverRange.Find.Execute(". However, ", $true, $false, $false, $false, $false, $true, 0, $false, ". Moreover, ", 2) | Out-Null

# --- re-find the now-updated text so we know exact character offsets -----
$ctxRange = $d.Content
$ctxRange.Find.Execute("6000. Moreover, upon", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if (-not $ctxRange.Find.Found) {
    throw "Could not locate the updated '6000. Moreover, upon' text"
}
$ctxStart = $ctxRange.Start

$dotSpaceStart  = $ctxStart + 4                                    # right after "6000"
$dotSpaceEnd    = $dotSpaceStart + 2                                # ". "
$moreoverStart  = $dotSpaceEnd
$moreoverEnd    = $moreoverStart + 8                                # "Moreover"
$commaSpaceStart= $moreoverEnd
$commaSpaceEnd  = $commaSpaceStart + 2                              # ", "

# --- locate the unique anchor text for edit #2 ("(Figure 2) ") -----------
$figRange = $d.Content
$figRange.Find.Execute("(Figure 2) ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if (-not $figRange.Find.Found) {
    throw "Could not locate the '(Figure 2) ' text"
}
$figStart = $figRange.Start
$figEnd   = $figRange.End

# ---------------------------------------------------------------------------
# Force run splits at the exact boundaries we need by toggling a formatting
# property on/off over a narrow Range (this splits the run without touching
# the text itself). Apply rightmost boundary first so earlier offsets stay
# valid.
# ---------------------------------------------------------------------------
function Split-RunBoundary($rangeStart, $rangeEnd) {
    $r = $d.Range($rangeStart, $rangeEnd)
    $r.Bold = 1
    $r.Bold = 0
}

Split-RunBoundary $figStart $figEnd
Split-RunBoundary $commaSpaceStart $commaSpaceEnd
Split-RunBoundary $moreoverStart $moreoverEnd
Split-RunBoundary $dotSpaceStart $dotSpaceEnd
